$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff represents a reordering (permutation) of the data rows (2-8),
# where each row's D (Fecha), M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado), Q (Unidad de comercialización), S (Precio $/Kg)
# and T (Kg / unidad) values move to a different row while all other columns
# stay the same (they were already identical across rows).

$rowData = @{
    2 = @{ D = 44418; M = 240; N = 10000; O = 11000; P = 10500; Q = "$/bandeja 10 kilos"; S = 1050; T = 10 }
    3 = @{ D = 44323; M = 270; N = 21000; O = 22000; P = 21500; Q = "$/bandeja 18 kilos"; S = 1194; T = 18 }
    4 = @{ D = 44487; M = 300; N = 14000; O = 15000; P = 14500; Q = "$/bandeja 10 kilos"; S = 1450; T = 10 }
    5 = @{ D = 44291; M = 200; N = 17000; O = 18000; P = 17500; Q = "$/bandeja 18 kilos"; S = 972;  T = 18 }
    6 = @{ D = 44263; M = 250; N = 21000; O = 22000; P = 21500; Q = "$/caja 18 kilos";    S = 1194; T = 18 }
    7 = @{ D = 44307; M = 250; N = 19000; O = 20000; P = 19500; Q = "$/bandeja 18 kilos"; S = 1083; T = 18 }
    8 = @{ D = 44489; M = 300; N = 26000; O = 27000; P = 26500; Q = "$/bandeja 18 kilos"; S = 1472; T = 18 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
